$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date strings in column A were re-formatted from DD/MM/YYYY to DD-MM-YYYY.
# These values must stay as plain text (they are not real Excel dates), so we
# force a text number format before assigning the new value, then restore the
# "Normal" style so no extra formatting is left behind on the cell.
$dates = @{
  3  = "28-07-2022"
  4  = "01-08-2022"
  5  = "04-08-2022"
  6  = "08-08-2022"
  7  = "11-08-2022"
  8  = "15-08-2022"
  9  = "18-08-2022"
  10 = "22-08-2022"
  11 = "25-08-2022"
  12 = "29-08-2022"
  13 = "01-09-2022"
  14 = "05-09-2022"
  15 = "08-09-2022"
  16 = "12-09-2022"
  17 = "15-09-2022"
  18 = "19-09-2022"
  19 = "22-09-2022"
  20 = "26-09-2022"
  21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.NumberFormat = "@"
  $cell.Value = $dates[$row]
  $cell.Style = "Normal"
}

# Update the attendance counters for rows 3-6 to match the refreshed data.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0
